# Reprocess rhyolitic glasses and update summary figure
# Updates the detection-limit summary rows (average / stdev / minimum / maximum,
# each stored at rows 2, 5 and 8) on every worksheet with the reprocessed values.

$wb = $excel.ActiveWorkbook

# Sheet 1: A870_2_bg_detlim
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 2).Value = 0.033
$ws.Cells.Item(2, 3).Value = 0.016
$ws.Cells.Item(2, 4).Value = 0.02
$ws.Cells.Item(2, 5).Value = 0.018
$ws.Cells.Item(2, 6).Value = 0.021
$ws.Cells.Item(2, 7).Value = 0.018
$ws.Cells.Item(2, 8).Value = 0.018
$ws.Cells.Item(2, 9).Value = 0.021
$ws.Cells.Item(2, 10).Value = 0.021
$ws.Cells.Item(2, 11).Value = 0.005
$ws.Cells.Item(2, 12).Value = 0.016
$ws.Cells.Item(2, 13).Value = 0.033
$ws.Cells.Item(5, 2).Value = 0.033
$ws.Cells.Item(5, 3).Value = 0.016
$ws.Cells.Item(5, 4).Value = 0.02
$ws.Cells.Item(5, 5).Value = 0.018
$ws.Cells.Item(5, 6).Value = 0.021
$ws.Cells.Item(5, 7).Value = 0.018
$ws.Cells.Item(5, 8).Value = 0.018
$ws.Cells.Item(5, 9).Value = 0.021
$ws.Cells.Item(5, 10).Value = 0.021
$ws.Cells.Item(5, 11).Value = 0.005
$ws.Cells.Item(5, 12).Value = 0.016
$ws.Cells.Item(5, 13).Value = 0.033
$ws.Cells.Item(8, 2).Value = 0.037
$ws.Cells.Item(8, 3).Value = 0.018
$ws.Cells.Item(8, 4).Value = 0.023
$ws.Cells.Item(8, 5).Value = 0.02
$ws.Cells.Item(8, 6).Value = 0.023
$ws.Cells.Item(8, 7).Value = 0.02
$ws.Cells.Item(8, 8).Value = 0.02
$ws.Cells.Item(8, 9).Value = 0.023
$ws.Cells.Item(8, 10).Value = 0.023
$ws.Cells.Item(8, 11).Value = 0.006
$ws.Cells.Item(8, 12).Value = 0.018
$ws.Cells.Item(8, 13).Value = 0.037

# Sheet 2: A870_3_bg_apf_detlim
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 2).Value = 0.04
$ws.Cells.Item(2, 3).Value = 0.02
$ws.Cells.Item(2, 4).Value = 0.024
$ws.Cells.Item(2, 5).Value = 0.022
$ws.Cells.Item(2, 6).Value = 0.025
$ws.Cells.Item(2, 7).Value = 0.022
$ws.Cells.Item(2, 8).Value = 0.021
$ws.Cells.Item(2, 9).Value = 0.025
$ws.Cells.Item(2, 10).Value = 0.025
$ws.Cells.Item(2, 11).Value = 0.006
$ws.Cells.Item(2, 12).Value = 0.02
$ws.Cells.Item(2, 13).Value = 0.04
$ws.Cells.Item(5, 2).Value = 0.04
$ws.Cells.Item(5, 3).Value = 0.02
$ws.Cells.Item(5, 4).Value = 0.024
$ws.Cells.Item(5, 5).Value = 0.022
$ws.Cells.Item(5, 6).Value = 0.025
$ws.Cells.Item(5, 7).Value = 0.022
$ws.Cells.Item(5, 8).Value = 0.021
$ws.Cells.Item(5, 9).Value = 0.025
$ws.Cells.Item(5, 10).Value = 0.025
$ws.Cells.Item(5, 11).Value = 0.006
$ws.Cells.Item(5, 12).Value = 0.02
$ws.Cells.Item(5, 13).Value = 0.04
$ws.Cells.Item(8, 2).Value = 0.045
$ws.Cells.Item(8, 3).Value = 0.022
$ws.Cells.Item(8, 4).Value = 0.028
$ws.Cells.Item(8, 5).Value = 0.024
$ws.Cells.Item(8, 6).Value = 0.028
$ws.Cells.Item(8, 7).Value = 0.024
$ws.Cells.Item(8, 8).Value = 0.024
$ws.Cells.Item(8, 9).Value = 0.028
$ws.Cells.Item(8, 10).Value = 0.028
$ws.Cells.Item(8, 11).Value = 0.007
$ws.Cells.Item(8, 12).Value = 0.022
$ws.Cells.Item(8, 13).Value = 0.045

# Sheet 3: A876_2_bg_detlim
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 2).Value = 0.017
$ws.Cells.Item(2, 3).Value = 0.032
$ws.Cells.Item(2, 4).Value = 0.016
$ws.Cells.Item(2, 5).Value = 0.018
$ws.Cells.Item(2, 6).Value = 0.016
$ws.Cells.Item(2, 7).Value = 0.029
$ws.Cells.Item(2, 8).Value = 0.021
$ws.Cells.Item(2, 9).Value = 0.007
$ws.Cells.Item(2, 10).Value = 0.016
$ws.Cells.Item(2, 11).Value = 0.032
$ws.Cells.Item(5, 2).Value = 0.017
$ws.Cells.Item(5, 3).Value = 0.032
$ws.Cells.Item(5, 4).Value = 0.016
$ws.Cells.Item(5, 5).Value = 0.018
$ws.Cells.Item(5, 6).Value = 0.016
$ws.Cells.Item(5, 7).Value = 0.029
$ws.Cells.Item(5, 8).Value = 0.021
$ws.Cells.Item(5, 9).Value = 0.007
$ws.Cells.Item(5, 10).Value = 0.016
$ws.Cells.Item(5, 11).Value = 0.032
$ws.Cells.Item(8, 2).Value = 0.019
$ws.Cells.Item(8, 3).Value = 0.036
$ws.Cells.Item(8, 4).Value = 0.018
$ws.Cells.Item(8, 5).Value = 0.02
$ws.Cells.Item(8, 6).Value = 0.017
$ws.Cells.Item(8, 7).Value = 0.032
$ws.Cells.Item(8, 8).Value = 0.024
$ws.Cells.Item(8, 9).Value = 0.008
$ws.Cells.Item(8, 10).Value = 0.017
$ws.Cells.Item(8, 11).Value = 0.036

# Sheet 4: A876_3_bg_apf_detlim
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 2).Value = 0.021
$ws.Cells.Item(2, 3).Value = 0.039
$ws.Cells.Item(2, 4).Value = 0.019
$ws.Cells.Item(2, 5).Value = 0.021
$ws.Cells.Item(2, 6).Value = 0.019
$ws.Cells.Item(2, 7).Value = 0.035
$ws.Cells.Item(2, 8).Value = 0.026
$ws.Cells.Item(2, 9).Value = 0.008999999999999999
$ws.Cells.Item(2, 10).Value = 0.019
$ws.Cells.Item(2, 11).Value = 0.039
$ws.Cells.Item(5, 2).Value = 0.021
$ws.Cells.Item(5, 3).Value = 0.039
$ws.Cells.Item(5, 4).Value = 0.019
$ws.Cells.Item(5, 5).Value = 0.021
$ws.Cells.Item(5, 6).Value = 0.019
$ws.Cells.Item(5, 7).Value = 0.035
$ws.Cells.Item(5, 8).Value = 0.026
$ws.Cells.Item(5, 9).Value = 0.008999999999999999
$ws.Cells.Item(5, 10).Value = 0.019
$ws.Cells.Item(5, 11).Value = 0.039
$ws.Cells.Item(8, 2).Value = 0.023
$ws.Cells.Item(8, 3).Value = 0.043
$ws.Cells.Item(8, 4).Value = 0.021
$ws.Cells.Item(8, 5).Value = 0.024
$ws.Cells.Item(8, 6).Value = 0.021
$ws.Cells.Item(8, 7).Value = 0.039
$ws.Cells.Item(8, 8).Value = 0.028
$ws.Cells.Item(8, 9).Value = 0.01
$ws.Cells.Item(8, 10).Value = 0.021
$ws.Cells.Item(8, 11).Value = 0.043

# Sheet 5: B989_2_bg_detlim
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 2).Value = 0.017
$ws.Cells.Item(2, 3).Value = 0.031
$ws.Cells.Item(2, 4).Value = 0.029
$ws.Cells.Item(2, 5).Value = 0.018
$ws.Cells.Item(2, 6).Value = 0.017
$ws.Cells.Item(2, 7).Value = 0.023
$ws.Cells.Item(2, 8).Value = 0.023
$ws.Cells.Item(2, 9).Value = 0.006
$ws.Cells.Item(2, 10).Value = 0.017
$ws.Cells.Item(2, 11).Value = 0.031
$ws.Cells.Item(5, 2).Value = 0.017
$ws.Cells.Item(5, 3).Value = 0.031
$ws.Cells.Item(5, 4).Value = 0.029
$ws.Cells.Item(5, 5).Value = 0.018
$ws.Cells.Item(5, 6).Value = 0.017
$ws.Cells.Item(5, 7).Value = 0.023
$ws.Cells.Item(5, 8).Value = 0.023
$ws.Cells.Item(5, 9).Value = 0.006
$ws.Cells.Item(5, 10).Value = 0.017
$ws.Cells.Item(5, 11).Value = 0.031
$ws.Cells.Item(8, 2).Value = 0.021
$ws.Cells.Item(8, 3).Value = 0.038
$ws.Cells.Item(8, 4).Value = 0.035
$ws.Cells.Item(8, 5).Value = 0.022
$ws.Cells.Item(8, 6).Value = 0.021
$ws.Cells.Item(8, 7).Value = 0.028
$ws.Cells.Item(8, 8).Value = 0.027
$ws.Cells.Item(8, 9).Value = 0.008
$ws.Cells.Item(8, 10).Value = 0.021
$ws.Cells.Item(8, 11).Value = 0.038

# Sheet 6: B989_3_bg_apf_detlim
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 2).Value = 0.021
$ws.Cells.Item(2, 3).Value = 0.038
$ws.Cells.Item(2, 4).Value = 0.035
$ws.Cells.Item(2, 5).Value = 0.022
$ws.Cells.Item(2, 6).Value = 0.021
$ws.Cells.Item(2, 7).Value = 0.028
$ws.Cells.Item(2, 8).Value = 0.027
$ws.Cells.Item(2, 9).Value = 0.008
$ws.Cells.Item(2, 10).Value = 0.021
$ws.Cells.Item(2, 11).Value = 0.038
$ws.Cells.Item(5, 2).Value = 0.021
$ws.Cells.Item(5, 3).Value = 0.038
$ws.Cells.Item(5, 4).Value = 0.035
$ws.Cells.Item(5, 5).Value = 0.022
$ws.Cells.Item(5, 6).Value = 0.021
$ws.Cells.Item(5, 7).Value = 0.028
$ws.Cells.Item(5, 8).Value = 0.027
$ws.Cells.Item(5, 9).Value = 0.008
$ws.Cells.Item(5, 10).Value = 0.021
$ws.Cells.Item(5, 11).Value = 0.038
$ws.Cells.Item(8, 2).Value = 0.025
$ws.Cells.Item(8, 3).Value = 0.046
$ws.Cells.Item(8, 4).Value = 0.042
$ws.Cells.Item(8, 5).Value = 0.026
$ws.Cells.Item(8, 6).Value = 0.025
$ws.Cells.Item(8, 7).Value = 0.034
$ws.Cells.Item(8, 8).Value = 0.033
$ws.Cells.Item(8, 9).Value = 0.008999999999999999
$ws.Cells.Item(8, 10).Value = 0.025
$ws.Cells.Item(8, 11).Value = 0.046

